$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above old row 17 (CC / 1052730461 / SARITA / 2506 row).
# This pushes the old row 17 down to row 18 (keeping its original formatting),
# and pushes the footer rows (22,23) down to (23,24) as well.
$ws.Rows.Item(17).Insert()

# Copy the formatting (borders/number format/font) of row 16 into the new row 17
# so the new row matches the "middle of table" look instead of the old
# "last row" look that the inserted blank row inherited.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)  # xlPasteFormats

# New period row (2506) - same worker / salary data as the existing rows.
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1052730461"
$ws.Range("D17").Value = "SARITA EUGENIA CASTRO MURILLO"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# The previously existing "2505" row becomes the most-recent period "2507".
$ws.Range("E16").Value = "2507"

# The row that got pushed down to 18 (used to be row 17, period 2506) now
# represents the oldest period shown, "2505".
$ws.Range("E18").Value = "2505"

# Update the totals to reflect 3 periods of mora instead of 2.
$ws.Range("E11").Value = 170820
$ws.Range("F13").Value = 3

$wb.Save()
